$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- FarmerStatTable -> FarmerLevelTable migration & table cleanup ---
#
# The "level" column (old column D) becomes the second column (new column B),
# right after "id" (column A). The former "eggStoreLimit"/"farmerStoreLimit"
# columns (old B/C) shift right into new C/D. "gold"/"gem" (E/F) are untouched.
#
# Achieve this with a classic insert-cut-delete column shuffle:
#   1. Insert a blank column before column B (B..F -> C..G).
#   2. Cut the "level" data (now in column E) into the new blank column B.
#   3. Delete the now-empty column E, closing the gap back to a 6-column table.

$ws.Columns("B").Insert(-4121)                       # xlShiftToRight
$ws.Range("E2:E17").Cut($ws.Range("B2:B17"))
$ws.Columns("E").Delete()

# --- Freeze panes at C4 (2 columns / 3 rows frozen) ---
$ws.Range("C4").Select()
$excel.ActiveWindow.FreezePanes = $true

# Leave the final selection as it was recorded in the authored workbook.
$ws.Range("J9").Select()
